$d = $word.ActiveDocument

# Anchor: the last paragraph in the body currently holds "BILAGA 1 - Fridlysta arter"
$anchor = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraStarts = @(0,0,0,0,0,0,0,0,0,0,0,0,0)

# --- Pass 1: insert every new paragraph with its plain text and style ---

# paragraph 1
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 1"
$p.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'
$paraStarts[0] = $p.Range.Start
$anchor = $p

# paragraph 2
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$paraStarts[1] = $p.Range.Start
$anchor = $p

# paragraph 3
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$paraStarts[2] = $p.Range.Start
$anchor = $p

# paragraph 4
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$paraStarts[3] = $p.Range.Start
$anchor = $p

# paragraph 5
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$paraStarts[4] = $p.Range.Start
$anchor = $p

# paragraph 6
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$paraStarts[5] = $p.Range.Start
$anchor = $p

# paragraph 7
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 2"
$p.Range.Text = 'Referenser - knärot'
$paraStarts[6] = $p.Range.Start
$anchor = $p

# paragraph 8
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$paraStarts[7] = $p.Range.Start
$anchor = $p

# paragraph 9
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$paraStarts[8] = $p.Range.Start
$anchor = $p

# paragraph 10
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$paraStarts[9] = $p.Range.Start
$anchor = $p

# paragraph 11
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$paraStarts[10] = $p.Range.Start
$anchor = $p

# paragraph 12
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$paraStarts[11] = $p.Range.Start
$anchor = $p

# paragraph 13
$anchor.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
$paraStarts[12] = $p.Range.Start
$anchor = $p

# --- Pass 2: go back and italicise the spans that need it (formatting only, no further inserts) ---

# paragraph 3 italics
$sub = $d.Range($paraStarts[2] + 34, $paraStarts[2] + 116)
$sub.Font.Italic = $true
$sub = $d.Range($paraStarts[2] + 278, $paraStarts[2] + 483)
$sub.Font.Italic = $true
$sub = $d.Range($paraStarts[2] + 490, $paraStarts[2] + 608)
$sub.Font.Italic = $true

# paragraph 4 italics
$sub = $d.Range($paraStarts[3] + 205, $paraStarts[3] + 1070)
$sub.Font.Italic = $true

# paragraph 8 italics
$sub = $d.Range($paraStarts[7] + 33, $paraStarts[7] + 113)
$sub.Font.Italic = $true

# paragraph 9 italics
$sub = $d.Range($paraStarts[8] + 62, $paraStarts[8] + 176)
$sub.Font.Italic = $true

# paragraph 10 italics
$sub = $d.Range($paraStarts[9] + 117, $paraStarts[9] + 207)
$sub.Font.Italic = $true

# paragraph 11 italics
$sub = $d.Range($paraStarts[10] + 54, $paraStarts[10] + 121)
$sub.Font.Italic = $true

# paragraph 12 italics
$sub = $d.Range($paraStarts[11] + 22, $paraStarts[11] + 57)
$sub.Font.Italic = $true

# paragraph 13 italics
$sub = $d.Range($paraStarts[12] + 25, $paraStarts[12] + 61)
$sub.Font.Italic = $true

# --- update the date in the "first page" header ---
$headers = $d.Sections.Item(1).Headers
$firstPageHeader = $headers.Item(2)
$firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Output "done"
